# Generate Report for Handback
# Updates the localization-status workbook to reflect that the handback
# has completed: sets Status to "Handed back: in sync with en-US",
# populates the "Latest Target File" / "Latest Handback File" columns
# (F/G) with hyperlinked file names, and stamps the
# "Latest Handback DateTime" column (H) with the actual handback time.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

function Update-HandbackSheet($ws, $h2Time, $h3Time, $mdUrl1, $mdUrl2, $xlf1Url, $xlf2Url, $mdDisplay1, $mdDisplay2, $xlf1Display, $xlf2Display) {

    # Status column (C) -- now handed back
    $ws.Range("C2").Value = $statusText
    $ws.Range("C3").Value = $statusText

    # Latest Target File (F) -- same source markdown file that was handed off
    $ws.Hyperlinks.Add($ws.Range("F2"), $mdUrl1, "", "", $mdDisplay1)
    $ws.Hyperlinks.Add($ws.Range("F3"), $mdUrl2, "", "", $mdDisplay2)

    # Latest Handback File (G) -- the translated xlf file handed back
    $ws.Hyperlinks.Add($ws.Range("G2"), $xlf1Url, "", "", $xlf1Display)
    $ws.Hyperlinks.Add($ws.Range("G3"), $xlf2Url, "", "", $xlf2Display)

    # Latest Handback DateTime (H)
    $ws.Range("H2").Value = $h2Time
    $ws.Range("H3").Value = $h3Time
}

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
Update-HandbackSheet $wsZh "2016-03-11 10:41:48" "2016-03-11 10:41:48" `
    "https://github.com/OpenLocalizationTest/oltest/blob/abde6c3c5ebcd5f40fca1ada00179258d7144607/e2e/0b21aba7-7f3a-424d-92eb-a1266b37238e.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/abde6c3c5ebcd5f40fca1ada00179258d7144607/e2e/1f1040ba-0717-43ae-a2e3-60273562994b.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d4c496ce50d3258bb755fc47703ea040004d129e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0b21aba7-7f3a-424d-92eb-a1266b37238e.e66bd2204af96e7035b7422a7d9faa370150cab0.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d4c496ce50d3258bb755fc47703ea040004d129e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/1f1040ba-0717-43ae-a2e3-60273562994b.a974af46748ef54968655947d4fbc530e0071fe8.zh-cn.xlf" `
    "0b21aba7-7f3a-424d-92eb-a1266b37238e.md" `
    "1f1040ba-0717-43ae-a2e3-60273562994b.md" `
    "0b21aba7-7f3a-424d-92eb-a1266b37238e.e66bd2204af96e7035b7422a7d9faa370150cab0.zh-cn.xlf" `
    "1f1040ba-0717-43ae-a2e3-60273562994b.a974af46748ef54968655947d4fbc530e0071fe8.zh-cn.xlf"

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
Update-HandbackSheet $wsDe "2016-03-11 10:41:53" "2016-03-11 10:41:53" `
    "https://github.com/OpenLocalizationTest/oltest/blob/abde6c3c5ebcd5f40fca1ada00179258d7144607/e2e/0b21aba7-7f3a-424d-92eb-a1266b37238e.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/abde6c3c5ebcd5f40fca1ada00179258d7144607/e2e/1f1040ba-0717-43ae-a2e3-60273562994b.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dc28d683daed4edf8193c1a276cf3da0b4e93643/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0b21aba7-7f3a-424d-92eb-a1266b37238e.e66bd2204af96e7035b7422a7d9faa370150cab0.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dc28d683daed4edf8193c1a276cf3da0b4e93643/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/1f1040ba-0717-43ae-a2e3-60273562994b.a974af46748ef54968655947d4fbc530e0071fe8.de-de.xlf" `
    "0b21aba7-7f3a-424d-92eb-a1266b37238e.md" `
    "1f1040ba-0717-43ae-a2e3-60273562994b.md" `
    "0b21aba7-7f3a-424d-92eb-a1266b37238e.e66bd2204af96e7035b7422a7d9faa370150cab0.de-de.xlf" `
    "1f1040ba-0717-43ae-a2e3-60273562994b.a974af46748ef54968655947d4fbc530e0071fe8.de-de.xlf"

$wsZh.Select()

Write-Host "Handback report generated."
